# edit.ps1 - applies the commit "Get ready for html tweaking" changes
# to MatthewDavidLitwin.docx via Word COM interop.

$d = $word.ActiveDocument
$enDash = [char]0x2013
$trademark = [char]0x2122

# ---------------------------------------------------------------------
# 1) Insert a new "FirstParagraph" styled date line right after the
#    "McGraw Hill Education" heading: "07/31/2012 - Present"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("McGraw Hill Education", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pHeading = $r.Paragraphs(1)
$pHeading.Range.InsertParagraphAfter()
$pDate = $pHeading.Next()
$pDate.Style = "FirstParagraph"
$pDate.Range.Text = "07/31/2012 " + $enDash + " Present"

# ---------------------------------------------------------------------
# 2) Insert a new "FirstParagraph" styled, italic line right after the
#    "Software Engineering Manager" Heading4: "Team Lead and Manager
#    Responsibilites"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Software Engineering Manager", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pMgrHeading = $r.Paragraphs(1)
$pMgrHeading.Range.InsertParagraphAfter()
$pTeamLead = $pMgrHeading.Next()
$pTeamLead.Style = "FirstParagraph"
$pTeamLead.Range.Text = "Team Lead and Manager Responsibilites"
$ilStart = $pTeamLead.Range.Start
$ilEnd = $ilStart + ("Team Lead and Manager Responsibilites").Length
$d.Range($ilStart, $ilEnd).Font.Italic = 1

# ---------------------------------------------------------------------
# 3) "Lead for a diverse" -> "Lead a diverse"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Lead for a diverse", $true, $false, $false, $false, $false, $true, 1, $false, "Lead a diverse", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Drop the trailing period after "...infrastructure and support."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Backend, front-end, infrastructure and support.", $true, $false, $false, $false, $false, $true, 1, $false, "Backend, front-end, infrastructure and support", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Drop the trailing period after "...contractors, and vendors."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Management of full time employees, contractors, and vendors.", $true, $false, $false, $false, $false, $true, 1, $false, "Management of full time employees, contractors, and vendors", 2) | Out-Null

# ---------------------------------------------------------------------
# 6) "implementation: coding, and directing developers" -> capitalized
# ---------------------------------------------------------------------
$d.Content.Find.Execute("implementation: coding, and directing developers", $true, $false, $false, $false, $false, $true, 1, $false, "Implementation: coding, and directing developers", 2) | Out-Null

# ---------------------------------------------------------------------
# 7) "negotiating feature-set and schedule with stakeholders" -> capitalized
# ---------------------------------------------------------------------
$d.Content.Find.Execute("negotiating feature-set and schedule with stakeholders", $true, $false, $false, $false, $false, $true, 1, $false, "Negotiating feature-set and schedule with stakeholders", 2) | Out-Null

# ---------------------------------------------------------------------
# 8) "03/22/1999 to 07/31/2012" -> "03/22/1999 - 07/31/2012" (en dash)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("03/22/1999 to 07/31/2012", $true, $false, $false, $false, $false, $true, 1, $false, "03/22/1999 " + $enDash + " 07/31/2012", 2) | Out-Null

# ---------------------------------------------------------------------
# 9) "Junior programmer to Senior Software Engineer" -> capitalized
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Junior programmer to Senior Software Engineer", $true, $false, $false, $false, $false, $true, 1, $false, "Junior Programmer to Senior Software Engineer", 2) | Out-Null

# ---------------------------------------------------------------------
# 10) Italicize "The Geometer's Sketchpad" and append the description
#     run right after it.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("The Geometer's Sketchpad", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Italic = 1
$r.InsertAfter(", a math education Dynamic Geometery visualization and interaction application.")

# ---------------------------------------------------------------------
# 11) Italicize "Fathom (TM) Dynamic Data (TM) Software"
# ---------------------------------------------------------------------
$r = $d.Content
$fathom = "Fathom" + $trademark + " Dynamic Data" + $trademark + " Software"
$r.Find.Execute($fathom, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Font.Italic = 1

# ---------------------------------------------------------------------
# 12) "...with an in-house cross-platform framework" -> "...C++ framework"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("with an in-house cross-platform framework", $true, $false, $false, $false, $false, $true, 1, $false, "with an in-house cross-platform C++ framework", 2) | Out-Null

# ---------------------------------------------------------------------
# 13) Rename the "Contact and Sites" heading/bookmark to "Contacts"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Contact and Sites", $true, $false, $false, $false, $false, $true, 1, $false, "Contacts", 2) | Out-Null
$oldBookmark = $d.Bookmarks("contact-and-sites")
$bmRange = $oldBookmark.Range
$d.Bookmarks.Add("contacts", $bmRange)
$d.Bookmarks("contact-and-sites").Delete()

# ---------------------------------------------------------------------
# 14) Remove the trailing " |" runs after the "linkedin" hyperlink.
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("|", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$toDelete = $d.Range($r.Start - 1, $r.End)
$toDelete.Delete()
